$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 becomes the data that used to be in row 3 (Segunda / Provincia de Curicó)
$ws.Range("D2").Value = 45086
$ws.Range("L2").Value = "Segunda"
$ws.Range("N2").Value = 20000
$ws.Range("O2").Value = 21000
$ws.Range("P2").Value = 20500
$ws.Range("R2").Value = "Provincia de Curicó"
$ws.Range("S2").Value = 1139

# Row 3 becomes the data that used to be in row 2 (Primera / Región de O'Higgins)
$ws.Range("D3").Value = 45043
$ws.Range("L3").Value = "Primera"
$ws.Range("N3").Value = 19000
$ws.Range("O3").Value = 20000
$ws.Range("P3").Value = 19500
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 1083
